$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 350, shifting existing rows 350:370 down to 351:371
$ws.Rows.Item(350).Insert()

# Match the date-formatted style used by the neighboring "Fecha" cells
$ws.Range("D350").NumberFormat = $ws.Range("D351").NumberFormat

# Populate the newly inserted row 350 with the new weekly data point
$ws.Range("A350").Value = 10
$ws.Range("B350").Value = "Vega Modelo de Temuco"
$ws.Range("C350").Value = "La Araucanía"
$ws.Range("D350").Value = 45041
$ws.Range("E350").Value = 9
$ws.Range("F350").Value = 100112039
$ws.Range("G350").Value = "Ciboulette"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 20
$ws.Range("K350").Value = 5000
$ws.Range("L350").Value = 5000
$ws.Range("M350").Value = 5000
$ws.Range("N350").Value = "$/docena de atados"
$ws.Range("O350").Value = "Provincia de Cautín"
$ws.Range("P350").Value = 1667
$ws.Range("Q350").Value = 3
$ws.Range("R350").Value = "Hortaliza"
